$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D8 value (1280 -> 1200); dependent formulas in H8 and K8 recalc automatically.
$ws.Range("D8").Value = 1200

# Update the active selection on the sheet from E10 to G10.
$ws.Range("G10").Select()
